$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 125, pushing the existing rows 125-139 down to 126-140.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly record.
$ws.Range("A125").Value = 6
$ws.Range("B125").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C125").Value = "Metropolitana"
$ws.Range("D125").Value = 45131
$ws.Range("E125").Value = 13
$ws.Range("F125").Value = 100114007
$ws.Range("G125").Value = "Jengibre"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 440
$ws.Range("K125").Value = 14000
$ws.Range("L125").Value = 15000
$ws.Range("M125").Value = 14523
$ws.Range("N125").Value = "$/caja 13 kilos"
$ws.Range("O125").Value = "Perú"
$ws.Range("P125").Value = 1117
$ws.Range("Q125").Value = 13
$ws.Range("R125").Value = "Hortaliza"
